# Swap the full contents of row 5 and row 6 on the active sheet.
# Row 5 held the "Garnlav" observation (Id 131067826) and row 6 held the
# "Tretåig hackspett" observation (Id 131067473); the edit swaps their
# positions (all columns A:AY) so row 5 now holds the woodpecker record
# and row 6 holds the lichen record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Columns whose value differs between row 5 and row 6 and where both
#     sides are "real" (non-blank) values -> plain swap via Value2. ---
$simpleCols = @("A","B","E","F","G","H","P","Q","R","AW","AX")

foreach ($col in $simpleCols) {
    $c5 = $ws.Range($col + "5")
    $c6 = $ws.Range($col + "6")
    $v5 = $c5.Value2
    $v6 = $c6.Value2
    $c5.Value2 = $v6
    $c6.Value2 = $v5
}

# --- Columns where one row has real text and the other is blank. The
#     blank side must end up empty (clears to a true blank cell), the
#     other side receives the text. ---
$ws.Range("M5").Value2  = "färska spår"
$ws.Range("M6").Value2  = ""

$ws.Range("Z5").Value2  = "15:46"
$ws.Range("Z6").Value2  = ""

$ws.Range("AB5").Value2 = "15:46"
$ws.Range("AB6").Value2 = ""

$ws.Range("AC5").Value2 = "Färska ringhack"
$ws.Range("AC6").Value2 = ""

$ws.Range("AH5").Value2 = ""
$ws.Range("AH6").Value2 = "Granskog"

# --- Columns where one row has an explicit empty-string cell and the
#     other row has no cell at all (true blank). The blank side needs to
#     become an explicit empty string (leading apostrophe forces Excel to
#     store a literal, empty text value instead of leaving the cell
#     blank); the apostrophe also flags the cell as "quote prefixed", so
#     the style is reset back to Normal right after. The other side of
#     each pair is simply cleared to blank. ---
$ws.Range("J5").Value2  = ""

$ws.Range("J6").Value2  = "'"
$ws.Range("J6").Style   = "Normal"

$ws.Range("L5").Value2  = "'"
$ws.Range("L5").Style   = "Normal"

$ws.Range("L6").Value2  = ""

$ws.Range("AF5").Value2 = ""

$ws.Range("AF6").Value2 = "'"
$ws.Range("AF6").Style  = "Normal"
